$wb = $excel.ActiveWorkbook

# --- "info" sheet: update labels/values and add a new "description" row ---
$info = $wb.Worksheets.Item("info")

# Row 2: quantity -> "carbon footprint" (was "global warming potential")
$info.Cells.Item(2, 2).Value = "carbon footprint"

# Row 4: title -> "everyday climate impact" (was "Carbon Footprints")
$info.Cells.Item(4, 2).Value = "everyday climate impact"

# Row 6 (new): description -> long descriptive blurb
$info.Cells.Item(6, 1).Value = "description"
$info.Cells.Item(6, 2).Value = 'Ever wondered about the climate impact of your daily choices? This collection measures the [carbon footprint](https://en.wikipedia.org/wiki/Carbon_footprint) of various items and activities, such as food production, electricity usage, and transportation. It uses kilograms of \( CO_2 \) equivalent (\(kg CO_2eq\)) as a standard unit, allowing us to compare the warming effect of different greenhouse gases, like methane, to that of carbon dioxide. By comparing these values, you''ll learn which activities have a large impact on our planet''s climate and which are less important.'

# Row 5: relabel "description" -> "tagline" (value - the old description text - stays put)
$info.Cells.Item(5, 1).Value = "tagline"

# Wrap text + widen column B so the long text is readable; size column A to fit labels.
$info.Columns.Item(2).WrapText = $true
$info.Columns.Item(1).ColumnWidth = 13.1
$info.Columns.Item(2).ColumnWidth = 75.9
$info.Rows.Item(6).RowHeight = 102

# --- view-state: move selection on the other sheets, then land on "info" ---
$data = $wb.Worksheets.Item("data")
$data.Range("C25").Select()

$l10n = $wb.Worksheets.Item("l10n")
$l10n.Range("B20").Select()

$info.Activate()
$info.Range("B6").Select()
